# Updated symbol list on Sat Jan 21 23:30:28 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for the
# crypto rows on Sheet1 with newly scraped values. Cells are stored as plain
# text (e.g. "301.96", "-0.74%"), so NumberFormat is forced to "Text" (@)
# before assigning the value to stop Excel's automatic number/percentage
# parsing, then restored to "General" to match the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: cell address, new text value
$updates = @(
    @("D2",  "301.96"),    @("E2",  "-0.74%"),
    @("D3",  "35.61"),     @("E3",  "5.44%"),
    @("D4",  "5.096"),     @("E4",  "-0.41%"),
    @("D5",  "0.07711"),   @("E5",  "-1.43%"),
    @("D6",  "2.187"),     @("E6",  "-9.55%"),
    @("D7",  "8.030"),     @("E7",  "0.10%"),
    @("D8",  "4.022"),     @("E8",  "2.73%"),
    @("D9",  "0.9284"),    @("E9",  "-1.01%"),
    @("D10", "0.09273"),   @("E10", "-6.39%"),
    @("D11", "0.1824"),    @("E11", "2.57%"),
    @("D12", "0.08535"),   @("E12", "-0.28%"),
    @("D13", "0.03653"),
                            @("E14", "0.26%"),
    @("D15", "0.001511"),  @("E15", "2.01%"),
    @("D16", "0.005658"),  @("E16", "-2.57%"),
    @("D17", "3.479"),
                            @("E18", "-0.12%"),
                            @("E19", "2.84%"),
    @("D20", "0.1326"),    @("E20", "-1.58%"),
    @("D21", "4.590"),     @("E21", "6.87%"),
    @("D22", "0.2248"),    @("E22", "7.56%"),
    @("D23", "0.04689"),   @("E23", "1.36%"),
    @("D24", "0.001239"),  @("E24", "1.62%"),
    @("D25", "0.004489"),  @("E25", "1.82%"),
    @("D26", "0.0001310"), @("E26", "1.28%"),
                            @("E27", "-20.27%"),
    @("D39", "0.01723"),   @("E39", "0.08%"),
    @("D40", "0.04691"),   @("E40", "-3.01%"),
    @("D41", "0.007952"),  @("E41", "1.55%"),
    @("D42", "0.1405"),    @("E42", "-0.52%"),
    @("D43", "0.007721"),  @("E43", "-21.27%"),
    @("D44", "0.002230"),  @("E44", "7.79%"),
    @("D45", "0.008979"),  @("E45", "-1.55%"),
    @("D46", "0.00006244"),@("E46", "2.13%"),
                            @("E47", "0.93%"),
    @("D48", "5.353"),     @("E48", "91.62%"),
    @("D49", "0.002710"),  @("E49", "35.79%"),
    @("D50", "0.00002115"),@("E50", "0.93%"),
    @("D51", "0.0002014"), @("E51", "0.93%")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $text = $u[1]
    $cell = $ws.Range($addr)
    # Force text storage so the numeric-/percent-looking string isn't
    # reinterpreted as a number, then restore the default "General" look.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}
